# task buscar producto completada
# Adds a new "producto" row (row 3) to the AgregarProducto sheet:
#   A3 -> hyperlink to the MercadoLibre search URL (same as A2)
#   B3 -> the product name found
#   C3 -> hyperlink (mailto:) to the contact e-mail (same as C2)
# and leaves the selection on B4, ready for the next entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: url / producto / correo -------------------------------------

# A3: same URL used in A2 - set the text first, then attach the hyperlink
# so the existing shared string is reused instead of creating a new one.
$ws.Range("A3").Value = "https://www.mercadolibre.com.co/"
$ws.Hyperlinks.Add($ws.Range("A3"), "https://www.mercadolibre.com.co/")

# B3: the product that was found
$ws.Range("B3").Value = "Cámara Fujifilm Instax Mini 12 Color Rosa"

# C3: same contact e-mail used in C2
$ws.Range("C3").Value = "paangudi3@gmail.com"
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:paangudi3@gmail.com")

# Move the active selection to B4, like after finishing data entry on row 3.
[void]$ws.Range("B4").Select()
